# Regenerate s_vals data to filter save games.
# Updates the numeric value columns (B:E and the computed sum column G)
# for data rows 2-9 on the active worksheet. Column F (Win) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => B, C, D, E, G (F stays as-is)
$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    3 = @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 13.86384647080068, 17.50910633199374)
    4 = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 13.86384647080068, 19.36876847130326)
    5 = @(0.2881169905109251, 1.626987699542094, 3.223369029078222, 13.86384647080068, 19.00232018993193)
    6 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    7 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    8 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 19.48425592650926)
    9 = @(0.2881169905109251, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 1.012145535086602)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G
}
